$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row above the old row 17 (the "Note: Without an external
# clock..." row), pushing it down to row 18. This also shifts the
# mergeCell reference automatically.
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new checklist item.
$ws.Range("A14").Value = "External Oscilloscope"

# Restore the window view state described by the diff: scrolled down a
# bit (topLeftCell = A3) with the active cell now on the new row (A16).
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("A16").Select() | Out-Null
